$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column B (Taxonsorteringsordning) simply increments from 98930 to 98931.
$simpleRows = @(2,3,4,5,6,7,8,9,10,13,14)
foreach ($r in $simpleRows) {
    $ws.Cells.Item($r, 2).Value = 98931
}

# Rows 11 and 12 have their full data swapped (all columns except B),
# while column B gets new values: row 11 -> 98931, row 12 -> 79244.
$row11 = @{
    A = 130965861
    B = 98931
    D = "LC"
    E = 219790
    F = "Fläcknycklar"
    G = "Dactylorhiza maculata"
    H = "(L.) Soó"
    Q = 497138
    R = 6713448
    AC = "Betydelsefulla förekomster . inventering åt vasa vind"
    AX = "Anders Esplund, Pia Edfors, Enviro Planning"
}

$row12 = @{
    A = 130965935
    B = 79244
    D = "NT"
    E = 6425
    F = "Garnlav"
    G = "Alectoria sarmentosa"
    H = "(Ach.) Ach."
    Q = 496969
    R = 6713674
    AC = "Måttlig förekomst . inventering åt vasa vind"
    AX = "Pia Edfors, Enviro Planning"
}

$ws.Range("A11").Value = $row11.A
$ws.Range("B11").Value = $row11.B
$ws.Range("D11").Value = $row11.D
$ws.Range("E11").Value = $row11.E
$ws.Range("F11").Value = $row11.F
$ws.Range("G11").Value = $row11.G
$ws.Range("H11").Value = $row11.H
$ws.Range("Q11").Value = $row11.Q
$ws.Range("R11").Value = $row11.R
$ws.Range("AC11").Value = $row11.AC
$ws.Range("AX11").Value = $row11.AX

$ws.Range("A12").Value = $row12.A
$ws.Range("B12").Value = $row12.B
$ws.Range("D12").Value = $row12.D
$ws.Range("E12").Value = $row12.E
$ws.Range("F12").Value = $row12.F
$ws.Range("G12").Value = $row12.G
$ws.Range("H12").Value = $row12.H
$ws.Range("Q12").Value = $row12.Q
$ws.Range("R12").Value = $row12.R
$ws.Range("AC12").Value = $row12.AC
$ws.Range("AX12").Value = $row12.AX
